# Weekly update: add two new daily-price records (rows) for Coliflor,
# "Agrícola del Norte S.A. de Arica", pushing the prior records down by
# two rows (dimension grows from A1:R104 to A1:R106).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 89; this shifts
# the existing rows 89-104 down to 91-106, just like the diff shows.
$ws.Rows("89:90").Insert()

# --- New row 89 : Coliflor, Segunda ---
$ws.Cells.Item(89, 1).Value = 1
$ws.Cells.Item(89, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(89, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(89, 4).Value = 44644
$ws.Cells.Item(89, 5).Value = 15
$ws.Cells.Item(89, 6).Value = 100112008
$ws.Cells.Item(89, 7).Value = "Coliflor"
$ws.Cells.Item(89, 8).Value = "Sin especificar"
$ws.Cells.Item(89, 9).Value = "Segunda"
$ws.Cells.Item(89, 10).Value = 900
$ws.Cells.Item(89, 11).Value = 900
$ws.Cells.Item(89, 12).Value = 1000
$ws.Cells.Item(89, 13).Value = 950
$ws.Cells.Item(89, 14).Value = "$/unidad"
$ws.Cells.Item(89, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(89, 16).Value = 950
$ws.Cells.Item(89, 17).Value = 1
$ws.Cells.Item(89, 18).Value = "Hortaliza"

# --- New row 90 : Coliflor, Tercera ---
$ws.Cells.Item(90, 1).Value = 1
$ws.Cells.Item(90, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(90, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(90, 4).Value = 44644
$ws.Cells.Item(90, 5).Value = 15
$ws.Cells.Item(90, 6).Value = 100112008
$ws.Cells.Item(90, 7).Value = "Coliflor"
$ws.Cells.Item(90, 8).Value = "Sin especificar"
$ws.Cells.Item(90, 9).Value = "Tercera"
$ws.Cells.Item(90, 10).Value = 1000
$ws.Cells.Item(90, 11).Value = 500
$ws.Cells.Item(90, 12).Value = 600
$ws.Cells.Item(90, 13).Value = 550
$ws.Cells.Item(90, 14).Value = "$/unidad"
$ws.Cells.Item(90, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(90, 16).Value = 550
$ws.Cells.Item(90, 17).Value = 1
$ws.Cells.Item(90, 18).Value = "Hortaliza"
